$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 481
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $excel.WorksheetFunction.Ln($cell.Value2)
}
